# Working hours - add new bookings for the end of January 2024 and move
# the "Total" row down to account for the 6 new rows of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "Total" row is row 23. Insert 6 new rows above it so the
# new bookings can be appended before the total, and the total row ends
# up at row 29.
$ws.Rows.Item(23).Resize(6).Insert()

# Inserted rows don't inherit formatting automatically - copy it down from
# the row directly above (the last real data row) like Excel's normal
# "insert copied cells" behaviour would.
$ws.Rows.Item(22).Copy()
$ws.Range("A23:F28").PasteSpecial(-4122)

# New booking rows (dates are Excel serial values; times are fractions of
# a day matching the existing data in the sheet).
$newRows = @(
    @{ Row = 23; Date = 45316; From = 0.333333333333333; To = 0.666666666666667 },
    @{ Row = 24; Date = 45318; From = 0.333333333333333; To = 0.666666666666667 },
    @{ Row = 25; Date = 45319; From = 0.333333333333333; To = 0.666666666666667 },
    @{ Row = 26; Date = 45320; From = 0.416666666666667; To = 0.583333333333333 },
    @{ Row = 27; Date = 45321; From = 0.416666666666667; To = 0.583333333333333 },
    @{ Row = 28; Date = 45322; From = 0.416666666666667; To = 0.583333333333333 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.From
    $ws.Cells.Item($row, 3).Value = $r.To
    $ws.Cells.Item($row, 4).Formula = "=(C$row<B$row)+C$row-B$row"
    $ws.Cells.Item($row, 5).Value = 13.5
    $ws.Cells.Item($row, 6).Formula = "=(D$row*24)*E$row"
    # Assigning this formula (which multiplies a time-formatted cell)
    # otherwise drags D's [hh]:mm:ss number format onto F - put it back to
    # the plain "General" style that the rest of column F uses.
    $ws.Cells.Item($row, 6).NumberFormat = "General"
}

# The total row is now row 29; fix up the SUM ranges to include the newly
# inserted rows (2..28 instead of 2..20).
$ws.Cells.Item(29, 4).Formula = "=SUM(D2:D28)"
$ws.Cells.Item(29, 6).Formula = "=SUM(F2:F28)"

# Update the view to match the new selection/scroll position.
$ws.Cells.Item(30, 6).Select()
$excel.ActiveWindow.ScrollRow = 2
